$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C56").Value = 2476.33
$ws.Range("A58:E61").ClearContents()
